$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$pairs = @(
    ,@("2023-11-29 Wednesday", "2023-11-30 Thursday")
    ,@("11+35=46", "16+6=22")
    ,@("25-10=15", "81-26=55")
    ,@("38-33=5", "97-23=74")
    ,@("21+14=35", "68-12=56")
    ,@("58+20=78", "56-25=31")
    ,@("59-52=7", "79-45=34")
    ,@("48-16=32", "3+32=35")
    ,@("74-29=45", "45-26=19")
    ,@("98-1=97", "89-45=44")
    ,@("22+56=78", "93-89=4")
    ,@("69+19=88", "49+7=56")
    ,@("75-29=46", "62-55=7")
    ,@("98-54=44", "26+33=59")
    ,@("5+35=40", "71-60=11")
    ,@("34+13=47", "47-25=22")
    ,@("37-14=23", "64+32=96")
    ,@("57-6=51", "74-11=63")
    ,@("25-20=5", "92-54=38")
    ,@("54-41=13", "36-17=19")
    ,@("33-11=22", "37+59=96")
    ,@("37+34=71", "34+13=47")
    ,@("48+23=71", "70-45=25")
    ,@("87-79=8", "91-88=3")
    ,@("44+31=75", "41-16=25")
    ,@("72-64=8", "3+28=31")
    ,@("21+61=82", "56-8=48")
    ,@("3+87=90", "91-53=38")
    ,@("69-31=38", "43+21=64")
    ,@("58+16=74", "68+12=80")
    ,@("0+83=83", "71-13=58")
    ,@("2+28=30", "97-57=40")
    ,@("75-12=63", "9+36=45")
    ,@("12+79=91", "86-82=4")
    ,@("15+63=78", "84-78=6")
    ,@("54-35=19", "21+71=92")
    ,@("23+76=99", "30+24=54")
    ,@("46+49=95", "58+6=64")
    ,@("33+6=39", "25+39=64")
    ,@("39+22=61", "0+37=37")
    ,@("13-5=8", "71-50=21")
    ,@("80+7=87", "38+55=93")
    ,@("26-23=3", "65-27=38")
    ,@("45+52=97", "44+47=91")
    ,@("85-76=9", "43+38=81")
    ,@("50+28=78", "40+58=98")
    ,@("71+28=99", "95-84=11")
    ,@("0+64=64", "28-13=15")
    ,@("76-71=5", "9-9=0")
    ,@("32+36=68", "41+21=62")
    ,@("97-48=49", "11-8=3")
    ,@("0+63=63", "17+28=45")
    ,@("26-26=0", "37+4=41")
    ,@("36+16=52", "9+57=66")
    ,@("67-40=27", "43+48=91")
    ,@("10+80=90", "10+11=21")
    ,@("28+45=73", "89-60=29")
    ,@("94-90=4", "25+25=50")
    ,@("79-50=29", "73-43=30")
    ,@("95-76=19", "70-21=49")
    ,@("53-21=32", "69+25=94")
    ,@("60+21=81", "75-3=72")
    ,@("48-32=16", "43+40=83")
    ,@("47+34=81", "91-33=58")
    ,@("50+9=59", "23+58=81")
    ,@("0+56=56", "25+32=57")
    ,@("8+31=39", "1+95=96")
    ,@("49+36=85", "19-14=5")
    ,@("44+24=68", "63+20=83")
    ,@("52-30=22", "64-61=3")
    ,@("65-62=3", "31+46=77")
    ,@("34+38=72", "57+8=65")
    ,@("50+48=98", "74-42=32")
    ,@("81-4=77", "17+62=79")
    ,@("13+47=60", "54+43=97")
    ,@("29+24=53", "78+13=91")
    ,@("54+17=71", "55-41=14")
    ,@("46+18=64", "73+16=89")
    ,@("72-9=63", "51+17=68")
    ,@("55-0=55", "68-30=38")
    ,@("26+37=63", "44+49=93")
    ,@("61+11=72", "91+0=91")
    ,@("76-67=9", "27-18=9")
    ,@("75+9=84", "30+48=78")
    ,@("66-20=46", "32-31=1")
    ,@("57+38=95", "77-66=11")
    ,@("36+26=62", "50-5=45")
    ,@("65+30=95", "90-79=11")
    ,@("30+31=61", "78-62=16")
    ,@("40+8=48", "9+35=44")
    ,@("77-51=26", "25+15=40")
    ,@("37+28=65", "3+48=51")
    ,@("54-4=50", "94-9=85")
    ,@("61-1=60", "13-4=9")
    ,@("84-41=43", "86-24=62")
    ,@("2+72=74", "89-29=60")
    ,@("41+0=41", "17+70=87")
    ,@("75+15=90", "68-48=20")
    ,@("70+14=84", "5+29=34")
    ,@("61-10=51", "31+11=42")
    ,@("75-66=9", "6+65=71")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Host "Done applying $($pairs.Count) replacements"